$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.883.24"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "3.175.71"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "220.18"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").Value = "625.04"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").Value = "1.09"
$ws.Range("E7").Value = "  +24.27%  "
$ws.Range("D8").Value = "0.374"
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "3.173.86"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").Value = "0.768"
$ws.Range("E11").Value = "  +16.80%  "
$ws.Range("E12").Value = "  +5.85%  "
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +3.87%  "
$ws.Range("D14").Value = "35.40"
$ws.Range("E14").Value = "  +8.99%  "
$ws.Range("D15").Value = "5.59"
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("D16").Value = "90.416.96"
$ws.Range("D17").Value = "3.766.59"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "3.173.90"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "3.79"
$ws.Range("E19").Value = "  +11.34%  "
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "14.37"
$ws.Range("E21").Value = "  +6.12%  "
$ws.Range("D22").Value = "442.74"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "8.98"
$ws.Range("E23").Value = "  +8.77%  "
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("E25").Value = "  +9.85%  "
$ws.Range("D26").Value = "12.47"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("D27").Value = "87.11"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").Value = "3.351.23"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "9.34"
$ws.Range("E31").Value = "  +14.60%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -8.91%  "
$ws.Range("D33").Value = "526.87"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "25.05"
$ws.Range("E34").Value = "  +9.02%  "
$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D35").Value = "3.76"
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "7.06"
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.145"
$ws.Range("E37").Value = "  +9.77%  "
$ws.Range("E38").Value = "  +5.92%  "
$ws.Range("D39").Value = "1.31"
$ws.Range("E39").Value = "  +5.40%  "
$ws.Range("D40").Value = "0.174"
$ws.Range("E40").Value = "  +22.74%  "
$ws.Range("D41").Value = "22.21"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0854"
$ws.Range("E42").Value = "  +21.32%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "0.413"
$ws.Range("E44").Value = "  +11.28%  "
$ws.Range("D45").Value = "1.95"
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "149.16"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "1.35"
$ws.Range("E48").Value = "  +10.01%  "
$ws.Range("D49").Value = "44.20"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").Value = "4.37"
$ws.Range("E50").Value = "  +8.72%  "
$ws.Range("D51").Value = "0.650"
$ws.Range("E51").Value = "  +9.50%  "